# Applies the commit:
#   - 3 tables (slides 14, 15, 16) switch from tableStyleId
#     {2EDF45AC-EC04-4C5C-A400-83CCCFE4C80E} to {3E4481D4-A51E-486D-810C-8C2E1DE1BAD0}
#   - the two theme parts (ppt/theme/theme1.xml "Office Theme" and
#     ppt/theme/theme2.xml "Integral") swap their colour-scheme content.

$p = $ppt.ActivePresentation

# --- 1. Re-point the three tables at the new table style -------------------
$newStyleId = "{3E4481D4-A51E-486D-810C-8C2E1DE1BAD0}"
foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2. Swap the theme colour schemes between theme1.xml & theme2.xml ------
# The deck's active theme (ppt/theme/theme2.xml, "Integral") is reachable
# through Slide.ThemeColorScheme; its 12 slots are swapped to the colours
# that previously lived in theme1.xml ("Office Theme").
$officeColors = 0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477

$slide1 = $p.Slides.Item(1)
$tcs = $slide1.ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}
